{"js": "// Insert a \"Varsel: \" lead-in as its own run immediately before the\n// existing text of the paragraph that reads exactly\n// \"Kassasjoner er registrert.\" (the short, standalone paragraph near the\n// end of the document \u2014 not the earlier \"Ingen kassasjoner er\n// registrert.\" paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === \"Kassasjoner er registrert.\") {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find paragraph \"Kassasjoner er registrert.\"');\n}\n\n// Insert the new text at the very start of the paragraph.\nconst inserted = target.getRange(\"Start\").insertText(\"Varsel: \", Word.InsertLocation.before);\n\n// Force Word to keep the newly inserted text in its own run (distinct\n// from the pre-existing \"Kassasjoner er registrert.\" run) by touching a\n// direct-formatting property and then restoring it. Without this, text\n// inserted right next to existing text with identical formatting is\n// silently merged back into a single run.\ninserted.font.bold = true;\nawait context.sync();\ninserted.font.bold = false;\nawait context.sync();\n", "ps1": "# Insert a \"Varsel: \" lead-in as its own run immediately before the\n# existing text of the paragraph that reads exactly\n# \"Kassasjoner er registrert.\" (the short, standalone paragraph near the\n# end of the document -- not the earlier \"Ingen kassasjoner er\n# registrert.\" paragraph).\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"Kassasjoner er registrert.`r\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find paragraph 'Kassasjoner er registrert.'\"\n}\n\n# Collapse to the very start of the paragraph and insert the new text.\n$ins = $target.Range.Duplicate\n$ins.Collapse(1)              # wdCollapseStart\n$ins.InsertBefore(\"Varsel: \")\n\n# Force Word to keep the newly inserted text in its own run (distinct\n# from the pre-existing \"Kassasjoner er registrert.\" run) by touching a\n# direct-formatting property and then restoring it. Without this, text\n# inserted right next to existing text with identical formatting is\n# silently merged back into a single run.\n$ins.Font.Bold = 1\n$ins.Font.Bold = 0\n"}
